$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Fecha Modificó" in D3 and template placeholder in D4,
# matching the formatting of the neighboring columns (C3 / C4)
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "Fecha Modificó"

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "{{item.FechaModificacion}}"

$excel.CutCopyMode = $false

# Update the named range to include the new column D
$wb.Names.Item("CostosToma").RefersTo = "=CostosToma!`$A`$4:`$D`$5"

# Extend the conditional formatting range to include the new column D
$cf = $ws.Range("A4:C4").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A4:D4"))

# Move the active selection to H13
$ws.Range("H13").Select()
